$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - General donations: update Donations value
$ws.Range("B2").Value = 1285003.96

# Rows 3-5 get reordered: Victory Drones, 1000 Drones for Ukraine, Mobile Shower Units
$ws.Range("A3").Value = "Victory Drones"
$ws.Range("B3").Value = 354136.55

$ws.Range("A4").Value = "1000 Drones for Ukraine"
$ws.Range("B4").Value = 18348.93

$ws.Range("A5").Value = "Mobile Shower Units"
$ws.Range("B5").Value = 4067.7

# Row 6 - Flight to Recovery: update Donations value
$ws.Range("B6").Value = 1287.09

# Row 7 ("--") is unchanged

# Row 8 changes from Admin/0/85 to Veteranius/217.76/0
$ws.Range("A8").Value = "Veteranius"
$ws.Range("B8").Value = 217.76
$ws.Range("C8").Value = 0

# New rows 9-20
$ws.Range("A9").Value = "grants"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 571000

$ws.Range("A10").Value = "Admin"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 390.04

$ws.Range("A11").Value = "bank fees & service charges"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 45.38

$ws.Range("A12").Value = "drone purchases"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 264502.3

$ws.Range("A13").Value = "events participation expenses"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 745

$ws.Range("A14").Value = "supplies & materials"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 20

$ws.Range("A15").Value = "car purchases"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 6474.59

$ws.Range("A16").Value = "lodging"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 1947.06

$ws.Range("A17").Value = "maling and delivery"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 318

$ws.Range("A18").Value = "supplies and materials"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 8870.139999999999

$ws.Range("A19").Value = "legal fees"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 600

$ws.Range("A20").Value = "transportation and parking"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 1216.49
